$d = $word.ActiveDocument

# Locate the paragraph that contains the final "</div>" line preceding the
# two blank lines, the "Catatan: ..." note, and the trailing empty
# paragraph. We find it by its text so the script isn't dependent on a
# hard-coded paragraph index.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "</div>") {
        $target = $p
    }
}

# Select everything from the end of that paragraph's text through the end
# of the document (the two blank paragraphs, the "Catatan:" note paragraph
# and the final empty/bold paragraph) and remove it - this is what happens
# when a user puts the cursor right after "</div>" and presses Ctrl+Shift+End
# then Delete.
$tailRange = $d.Range($target.Range.End, $d.Content.End)
$tailRange.Delete()

# After the deletion, the former last paragraph's mark formatting (Arial,
# bold, lang=en-ID) must become the paragraph-mark formatting of the
# now-last paragraph (the "</div>" one), exactly like Word does when you
# merge a paragraph into the one before it by deleting the break. The run
# text/formatting itself ("</div>") is left completely untouched.
$count = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($count)

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' + `
    'w14:paraId="719B6F1E" w14:textId="5929FEE6" w:rsidR="00CB0DF1" w:rsidRDefault="000E4B7F" w:rsidP="000E4B7F" ' + `
    'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
    '<w:pPr><w:ind w:left="284" w:hanging="283"/>' + `
    '<w:rPr><w:rFonts w:ascii="Arial"/><w:b/><w:lang w:val="en-ID"/></w:rPr>' + `
    '</w:pPr>' + `
    '<w:r w:rsidRPr="000E4B7F"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="en-ID"/></w:rPr><w:t>&lt;/div&gt;</w:t></w:r>' + `
    '</w:p>'

$last.Range.InsertXML($newXml)
